# Auto-generated edit script: updates leve-profit calculation columns (H:N)
# across ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets, per the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 441.1
$ws.Range("I8").Value = 101.833336
$ws.Range("J8").Value = 950
$ws.Range("K8").Value = 305.500008
$ws.Range("L8").Value = 2850
$ws.Range("M8").Value = -166.500008
$ws.Range("N8").Value = -3128
$ws.Range("H21").Value = 9203.4
$ws.Range("I21").Value = 5008.5
$ws.Range("K21").Value = 5008.5
$ws.Range("M21").Value = -4540.5
$ws.Range("H23").Value = 9203.4
$ws.Range("I23").Value = 5008.5
$ws.Range("K23").Value = 5008.5
$ws.Range("M23").Value = -4774.5
$ws.Range("H40").Value = 3618.182
$ws.Range("I40").Value = 3850
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 3850
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -3675
$ws.Range("N40").Value = -3350
$ws.Range("H43").Value = 2315
$ws.Range("I43").Value = 3126.8333
$ws.Range("J43").Value = 1340.8
$ws.Range("K43").Value = 3126.8333
$ws.Range("L43").Value = 1340.8
$ws.Range("M43").Value = -3057.8333
$ws.Range("N43").Value = -1478.8
$ws.Range("H52").Value = 686979.3
$ws.Range("I52").Value = 686979.3
$ws.Range("K52").Value = 2060937.9
$ws.Range("M52").Value = -2060777.9
$ws.Range("H58").Value = 981946.1
$ws.Range("I58").Value = 2451378
$ws.Range("J58").Value = 2325
$ws.Range("K58").Value = 7354134
$ws.Range("L58").Value = 6975
$ws.Range("M58").Value = -7353984
$ws.Range("N58").Value = -7275
$ws.Range("H96").Value = 1252.7778
$ws.Range("I96").Value = 446.75
$ws.Range("J96").Value = 2864.8333
$ws.Range("K96").Value = 1340.25
$ws.Range("L96").Value = 8594.499899999999
$ws.Range("M96").Value = 32.75
$ws.Range("N96").Value = -11340.4999
$ws.Range("H132").Value = 7151069
$ws.Range("I132").Value = 7821125
$ws.Range("K132").Value = 23463375
$ws.Range("M132").Value = -23460845
$ws.Range("H133").Value = 51895
$ws.Range("J133").Value = 51895
$ws.Range("L133").Value = 51895
$ws.Range("N133").Value = -62015
$ws.Range("H141").Value = 5063.3335
$ws.Range("I141").Value = 6055.8
$ws.Range("J141").Value = 3822.75
$ws.Range("K141").Value = 18167.4
$ws.Range("L141").Value = 11468.25
$ws.Range("M141").Value = -12987.4
$ws.Range("N141").Value = -21828.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23168.545
$ws.Range("I32").Value = 6694.965
$ws.Range("K32").Value = 6694.965
$ws.Range("M32").Value = -6407.965
$ws.Range("H74").Value = 817.2105
$ws.Range("I74").Value = 746.71875
$ws.Range("J74").Value = 1193.1666
$ws.Range("K74").Value = 746.71875
$ws.Range("L74").Value = 1193.1666
$ws.Range("M74").Value = 127.28125
$ws.Range("N74").Value = -2941.1666
$ws.Range("H77").Value = 817.2105
$ws.Range("I77").Value = 746.71875
$ws.Range("J77").Value = 1193.1666
$ws.Range("K77").Value = 3733.59375
$ws.Range("L77").Value = 5965.833000000001
$ws.Range("M77").Value = 634.40625
$ws.Range("N77").Value = -14701.833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I107").Value = 333334000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 333334000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -333332080
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1228
$ws.Range("I122").Value = 1270.6666
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 3811.9998
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -1361.9998
$ws.Range("N122").Value = -8200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 125130744
$ws.Range("I80").Value = 166840160
$ws.Range("J80").Value = 2490
$ws.Range("K80").Value = 166840160
$ws.Range("L80").Value = 2490
$ws.Range("M80").Value = -166839162
$ws.Range("N80").Value = -4486
$ws.Range("H83").Value = 125130744
$ws.Range("I83").Value = 166840160
$ws.Range("J83").Value = 2490
$ws.Range("K83").Value = 834200800
$ws.Range("L83").Value = 12450
$ws.Range("M83").Value = -834195808
$ws.Range("N83").Value = -22434

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1992.7858
$ws.Range("J22").Value = 1685.1
$ws.Range("L22").Value = 1685.1
$ws.Range("N22").Value = -2275.1
$ws.Range("H27").Value = 1992.7858
$ws.Range("J27").Value = 1685.1
$ws.Range("L27").Value = 1685.1
$ws.Range("N27").Value = -1899.1
$ws.Range("H55").Value = 237349.92
$ws.Range("I55").Value = 474204.4
$ws.Range("J55").Value = 495.41666
$ws.Range("K55").Value = 474204.4
$ws.Range("L55").Value = 495.41666
$ws.Range("M55").Value = -474031.4
$ws.Range("N55").Value = -841.41666
$ws.Range("H69").Value = 35081.5
$ws.Range("J69").Value = 35081.5
$ws.Range("L69").Value = 35081.5
$ws.Range("N69").Value = -36703.5
$ws.Range("H72").Value = 35081.5
$ws.Range("J72").Value = 35081.5
$ws.Range("L72").Value = 105244.5
$ws.Range("N72").Value = -113356.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 400751.8
$ws.Range("I81").Value = 1000000
$ws.Range("J81").Value = 250939.75
$ws.Range("K81").Value = 2000000
$ws.Range("L81").Value = 501879.5
$ws.Range("M81").Value = -1998939
$ws.Range("N81").Value = -504001.5
$ws.Range("H84").Value = 400751.8
$ws.Range("I84").Value = 1000000
$ws.Range("J84").Value = 250939.75
$ws.Range("K84").Value = 10000000
$ws.Range("L84").Value = 2509397.5
$ws.Range("M84").Value = -9994696
$ws.Range("N84").Value = -2520005.5
$ws.Range("H113").Value = 588.86664
$ws.Range("I113").Value = 443.66666
$ws.Range("J113").Value = 806.6667
$ws.Range("K113").Value = 1330.99998
$ws.Range("L113").Value = 2420.0001
$ws.Range("M113").Value = 839.0000199999999
$ws.Range("N113").Value = -6760.0001
$ws.Range("H126").Value = 1351.4546
$ws.Range("I126").Value = 1511.2
$ws.Range("J126").Value = 1009.1429
$ws.Range("K126").Value = 4533.6
$ws.Range("L126").Value = 3027.4287
$ws.Range("M126").Value = -2063.6
$ws.Range("N126").Value = -7967.4287
$ws.Range("H132").Value = 4764.048
$ws.Range("I132").Value = 5679.4287
$ws.Range("J132").Value = 2933.2856
$ws.Range("K132").Value = 17038.2861
$ws.Range("L132").Value = 8799.856800000001
$ws.Range("M132").Value = -14508.2861
$ws.Range("N132").Value = -13859.8568
$ws.Range("H136").Value = 1412.0469
$ws.Range("I136").Value = 567.7308
$ws.Range("J136").Value = 1989.7368
$ws.Range("K136").Value = 1703.1924
$ws.Range("L136").Value = 5969.2104
$ws.Range("M136").Value = 846.8075999999999
$ws.Range("N136").Value = -11069.2104
